# Automatic update of files.
# Applies the row-content permutation for rows 59-62 and 86-95
# (species-occurrence records were re-ordered upstream; row numbers
# stay fixed but the per-row field values rotate among the rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59
$ws.Cells.Item(59,1).Value2 = 111896618
$ws.Cells.Item(59,2).Value2 = 90332
$ws.Cells.Item(59,5).Value2 = 4769
$ws.Cells.Item(59,6).Value2 = "Svavelriska"
$ws.Cells.Item(59,7).Value2 = "Lactarius scrobiculatus"
$ws.Cells.Item(59,8).Value2 = "(Scop.:Fr.) Fr."
$ws.Cells.Item(59,17).Value2 = 574291.9436535498
$ws.Cells.Item(59,18).Value2 = 6703488.13161656

# Row 60
$ws.Cells.Item(60,1).Value2 = 111896589
$ws.Cells.Item(60,2).Value2 = 101141
$ws.Cells.Item(60,4).Value2 = "LC"
$ws.Cells.Item(60,5).Value2 = 222002
$ws.Cells.Item(60,6).Value2 = "Underviol"
$ws.Cells.Item(60,7).Value2 = "Viola mirabilis"
$ws.Cells.Item(60,8).Value2 = "L."
$ws.Cells.Item(60,17).Value2 = 574319.8708033765
$ws.Cells.Item(60,18).Value2 = 6703355.482874738

# Row 61
$ws.Cells.Item(61,1).Value2 = 111896682
$ws.Cells.Item(61,2).Value2 = 88966
$ws.Cells.Item(61,4).Value2 = "NT"
$ws.Cells.Item(61,5).Value2 = 5754
$ws.Cells.Item(61,6).Value2 = "Gultoppig fingersvamp"
$ws.Cells.Item(61,7).Value2 = "Ramaria testaceoflava"
$ws.Cells.Item(61,8).Value2 = "(Bres.) Corner"
$ws.Cells.Item(61,17).Value2 = 574256.1442139128
$ws.Cells.Item(61,18).Value2 = 6703474.568800534

# Row 62
$ws.Cells.Item(62,1).Value2 = 111896629
$ws.Cells.Item(62,17).Value2 = 574410.1934356905
$ws.Cells.Item(62,18).Value2 = 6703385.462655744

# Row 86
$ws.Cells.Item(86,1).Value2 = 111886390
$ws.Cells.Item(86,2).Value2 = 98535
$ws.Cells.Item(86,5).Value2 = 222498
$ws.Cells.Item(86,6).Value2 = "Blåsippa"
$ws.Cells.Item(86,7).Value2 = "Hepatica nobilis"
$ws.Cells.Item(86,8).Value2 = "Schreb."
$ws.Cells.Item(86,16).Value2 = "Gropbackagruvorna (Gropbackagruvorna), Gstr"
$ws.Cells.Item(86,17).Value2 = 574320.5008898397
$ws.Cells.Item(86,18).Value2 = 6703541.511550271
$ws.Cells.Item(86,19).Value2 = 1
$ws.Cells.Item(86,26).Value2 = "11:12"
$ws.Cells.Item(86,28).Value2 = "11:12"
$ws.Cells.Item(86,49).Value2 = "Annelie Hilmerby"
$ws.Cells.Item(86,50).Value2 = "Annelie Hilmerby"

# Row 87
$ws.Cells.Item(87,1).Value2 = 111886198
$ws.Cells.Item(87,2).Value2 = 88966
$ws.Cells.Item(87,4).Value2 = "NT"
$ws.Cells.Item(87,5).Value2 = 5754
$ws.Cells.Item(87,6).Value2 = "Gultoppig fingersvamp"
$ws.Cells.Item(87,7).Value2 = "Ramaria testaceoflava"
$ws.Cells.Item(87,8).Value2 = "(Bres.) Corner"
$ws.Cells.Item(87,16).Value2 = "Gropbackagruvorna (Gropbackagruvorna), Gstr"
$ws.Cells.Item(87,17).Value2 = 574244.409384649
$ws.Cells.Item(87,18).Value2 = 6703468.407307444
$ws.Cells.Item(87,49).Value2 = "Patric Engfeldt"
$ws.Cells.Item(87,50).Value2 = "Patric Engfeldt"

# Row 88
$ws.Cells.Item(88,1).Value2 = 111886372
$ws.Cells.Item(88,2).Value2 = 90687
$ws.Cells.Item(88,5).Value2 = 5964
$ws.Cells.Item(88,6).Value2 = "Fjällig taggsvamp s.str."
$ws.Cells.Item(88,7).Value2 = "Sarcodon imbricatus s.str."
$ws.Cells.Item(88,8).Value2 = "(L.:Fr.) P.Karst."
$ws.Cells.Item(88,16).Value2 = "Gropbackagruvorna (Gropbackagruvorna), Gstr"
$ws.Cells.Item(88,17).Value2 = 574261.3270212604
$ws.Cells.Item(88,18).Value2 = 6703462.833304818
$ws.Cells.Item(88,49).Value2 = "Patric Engfeldt"
$ws.Cells.Item(88,50).Value2 = "Patric Engfeldt"

# Row 89
$ws.Cells.Item(89,1).Value2 = 111885842
$ws.Cells.Item(89,2).Value2 = 101141
$ws.Cells.Item(89,5).Value2 = 222002
$ws.Cells.Item(89,6).Value2 = "Underviol"
$ws.Cells.Item(89,7).Value2 = "Viola mirabilis"
$ws.Cells.Item(89,8).Value2 = "L."
$ws.Cells.Item(89,16).Value2 = "Gropbackagruvorna (Gropbackagruvorna), Gstr"
$ws.Cells.Item(89,17).Value2 = 574289.9628196132
$ws.Cells.Item(89,18).Value2 = 6703416.053151045
$ws.Cells.Item(89,49).Value2 = "Patric Engfeldt"
$ws.Cells.Item(89,50).Value2 = "Patric Engfeldt"

# Row 90
$ws.Cells.Item(90,1).Value2 = 111896595
$ws.Cells.Item(90,2).Value2 = 101703
$ws.Cells.Item(90,5).Value2 = 222412
$ws.Cells.Item(90,6).Value2 = "Tibast"
$ws.Cells.Item(90,7).Value2 = "Daphne mezereum"
$ws.Cells.Item(90,8).Value2 = "L."
$ws.Cells.Item(90,17).Value2 = 574247.1857444055
$ws.Cells.Item(90,18).Value2 = 6703405.307550027

# Row 91
$ws.Cells.Item(91,1).Value2 = 111896581
$ws.Cells.Item(91,2).Value2 = 99413
$ws.Cells.Item(91,5).Value2 = 221235
$ws.Cells.Item(91,6).Value2 = "Vårärt"
$ws.Cells.Item(91,7).Value2 = "Lathyrus vernus"
$ws.Cells.Item(91,8).Value2 = "(L.) Bernh."
$ws.Cells.Item(91,17).Value2 = 574333.7783005711
$ws.Cells.Item(91,18).Value2 = 6703424.353483723

# Row 92
$ws.Cells.Item(92,1).Value2 = 111896668
$ws.Cells.Item(92,2).Value2 = 85089
$ws.Cells.Item(92,5).Value2 = 3762
$ws.Cells.Item(92,6).Value2 = "Olivspindling"
$ws.Cells.Item(92,7).Value2 = "Cortinarius venetus"
$ws.Cells.Item(92,8).Value2 = "(Fr.:Fr.) Fr."
$ws.Cells.Item(92,11).Value2 = ""
$ws.Cells.Item(92,16).Value2 = "Kratte masugn, Gstr"
$ws.Cells.Item(92,17).Value2 = 574317.8180265825
$ws.Cells.Item(92,18).Value2 = 6703431.427106799
$ws.Cells.Item(92,19).Value2 = 25
$ws.Cells.Item(92,26).Value2 = "00:00"
$ws.Cells.Item(92,28).Value2 = "00:00"
$ws.Cells.Item(92,49).Value2 = "Philipp Weiss"
$ws.Cells.Item(92,50).Value2 = "Philipp Weiss"

# Row 93
$ws.Cells.Item(93,1).Value2 = 111896610
$ws.Cells.Item(93,2).Value2 = 90332
$ws.Cells.Item(93,4).Value2 = "LC"
$ws.Cells.Item(93,5).Value2 = 4769
$ws.Cells.Item(93,6).Value2 = "Svavelriska"
$ws.Cells.Item(93,7).Value2 = "Lactarius scrobiculatus"
$ws.Cells.Item(93,8).Value2 = "(Scop.:Fr.) Fr."
$ws.Cells.Item(93,11).Value2 = ""
$ws.Cells.Item(93,16).Value2 = "Kratte masugn, Gstr"
$ws.Cells.Item(93,17).Value2 = 574272.2591996479
$ws.Cells.Item(93,18).Value2 = 6703411.74269451
$ws.Cells.Item(93,49).Value2 = "Philipp Weiss"
$ws.Cells.Item(93,50).Value2 = "Philipp Weiss"

# Row 94
$ws.Cells.Item(94,1).Value2 = 111896696
$ws.Cells.Item(94,2).Value2 = 98535
$ws.Cells.Item(94,5).Value2 = 222498
$ws.Cells.Item(94,6).Value2 = "Blåsippa"
$ws.Cells.Item(94,7).Value2 = "Hepatica nobilis"
$ws.Cells.Item(94,8).Value2 = "Schreb."
$ws.Cells.Item(94,11).Value2 = ""
$ws.Cells.Item(94,16).Value2 = "Kratte masugn, Gstr"
$ws.Cells.Item(94,17).Value2 = 574272.5440735799
$ws.Cells.Item(94,18).Value2 = 6703373.755373025
$ws.Cells.Item(94,49).Value2 = "Philipp Weiss"
$ws.Cells.Item(94,50).Value2 = "Philipp Weiss"

# Row 95
$ws.Cells.Item(95,1).Value2 = 111896614
$ws.Cells.Item(95,2).Value2 = 90332
$ws.Cells.Item(95,5).Value2 = 4769
$ws.Cells.Item(95,6).Value2 = "Svavelriska"
$ws.Cells.Item(95,7).Value2 = "Lactarius scrobiculatus"
$ws.Cells.Item(95,8).Value2 = "(Scop.:Fr.) Fr."
$ws.Cells.Item(95,11).Value2 = ""
$ws.Cells.Item(95,16).Value2 = "Kratte masugn, Gstr"
$ws.Cells.Item(95,17).Value2 = 574228.885558943
$ws.Cells.Item(95,18).Value2 = 6703430.096512586
$ws.Cells.Item(95,49).Value2 = "Philipp Weiss"
$ws.Cells.Item(95,50).Value2 = "Philipp Weiss"

